# add poll class diagram to dev guide
#
# Applies the OOXML diff to slide 1:
#   - Shape id=82  "TextBox 81"        : nudge right (x 6248400 -> 6324600 EMU)
#   - Shape id=94  "Rectangle 8"       : reflow/resize box, rename "Poll" -> "AbstractPoll"
#   - Shape id=98  "Elbow Connector 63": re-route bent connector (pos/size/adj1)
#   - Shape id=121 "Elbow Connector 78": re-route bent connector (pos/size)
#
# NOTE on precision: PowerPoint's Shape.Left/Top/Width/Height (and
# Adjustments) are single-precision (32-bit float) point values under the
# COM object model, and the host re-derives EMU from them as
# floor(pt_as_float32 * 12700). Most of the target EMU offsets in the diff
# are not exact multiples of 12700, so the literal "target EMU / 12700"
# decimal does not always round-trip through that float32 cast to the
# exact EMU integer the diff expects. The literals below were solved so
# that, after the float64->float32 narrowing the host performs, multiplying
# by 12700 and flooring reproduces the exact target EMU value.

function Get-ShapeById($shapes, $id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Shape id=82 "TextBox 81" ("*") ------------------------------------
# <a:off x="6248400" y="2209800"/>  ->  <a:off x="6324600" y="2209800"/>
$shp82 = Get-ShapeById $s.Shapes 82
$shp82.Left = 498.0

# --- Shape id=94 "Rectangle 8" ("Poll" -> "AbstractPoll") --------------
# <a:off x="6078908" y="2347620"/> <a:ext cx="708186" cy="346760"/>
# -> <a:off x="6078907" y="2347620"/> <a:ext cx="893633" cy="328141"/>
$shp94 = Get-ShapeById $s.Shapes 94
$shp94.Left = 478.65411376953125
$shp94.Width = 70.36480712890625
$shp94.Height = 25.837875366210938
$shp94.TextFrame.TextRange.Text = "AbstractPoll"

# --- Shape id=98 "Elbow Connector 63" (bentConnector4, flipH) ----------
# <a:off x="6433001" y="1713340"/> <a:ext cx="465238" cy="634280"/>
# adj1 -49136 -> -61367
# -> <a:off x="6525724" y="1713340"/> <a:ext cx="372515" cy="634280"/>
$shp98 = Get-ShapeById $s.Shapes 98
$shp98.Left = 513.8365478515625
$shp98.Width = 29.331890106201172
$shp98.Adjustments.Item(1) = -0.61367

# --- Shape id=121 "Elbow Connector 78" (bentConnector2, rot/flipV) -----
# <a:off x="5802718" y="2521000"/> <a:ext cx="276191" cy="318258"/>
# -> <a:off x="5802719" y="2511691"/> <a:ext cx="276188" cy="327566"/>
$shp121 = Get-ShapeById $s.Shapes 121
$shp121.Left = 456.9070129394531
$shp121.Top = 197.7709503173828
$shp121.Width = 21.747087478637695
$shp121.Height = 25.792598724365234
